# Update automatico via Actualizar 02-13-2021 12-44-31
#
# The "availability" timestamp column (D) gets refreshed on each run:
# a brand-new timestamp is recorded and the previously-recorded
# timestamps cascade down to the next block of rows.
#
#   rows 2-15  (newest) -> 44240.53084095146
#   rows 16-29           -> 44240.5096403125
#   rows 30-43 (oldest)  -> 44240.48844734954

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value  = 44240.53084095146
$ws.Range("D16:D29").Value = 44240.5096403125
$ws.Range("D30:D43").Value = 44240.48844734954
